# Update the workbook "Översikt FORSHAGA" sheet:
#  1. Bump the "Förändrad" date (column C) from 2023-09-19 (45188) to
#     2023-09-20 (45189) for every existing data row (rows 2-300).
#  2. Give row 300 an explicit custom row height (matches the rest of
#     the data rows).
#  3. Append a brand-new data row (301) for case "A 44281-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the "changed" date for every existing data row in one shot.
$ws.Range("C2:C300").Value = 45189

# 2. Row 300 picks up an explicit custom height, like its neighbours.
$ws.Rows.Item(300).RowHeight = 15

# 3. Append the new row (301) with its values.
$ws.Range("A301").Value = "A 44281-2023"
$ws.Range("B301").Value = 45188
$ws.Range("C301").Value = 45189
$ws.Range("D301").Value = "VÄRMLANDS LÄN"
$ws.Range("E301").Value = "FORSHAGA"
$ws.Range("F301").Value = "Bergvik skog väst AB"
$ws.Range("G301").Value = 2.1
$ws.Range("H301").Value = 0
$ws.Range("I301").Value = 0
$ws.Range("J301").Value = 0
$ws.Range("K301").Value = 0
$ws.Range("L301").Value = 0
$ws.Range("M301").Value = 0
$ws.Range("N301").Value = 0
$ws.Range("O301").Value = 0
$ws.Range("P301").Value = 0
$ws.Range("Q301").Value = 0
$ws.Range("R301").Value = ""

# Match the date/number formatting used by the rest of the table.
$ws.Range("B301:C301").NumberFormat = "YYYY-MM-DD"
$ws.Range("R301").WrapText = $true
